# Update "想去人数" (interested-count) figures and a refreshed cover image
# URL across the workbook's sheets, matching the upstream data refresh
# ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 911
$ws1.Range("F7").Value = 464
$ws1.Range("F9").Value = 2191
$ws1.Range("F13").Value = 1088
$ws1.Range("F16").Value = 668
$ws1.Range("F17").Value = 13171
$ws1.Range("F18").Value = 1261
$ws1.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202409/UKoa3flf1726049294866.jpeg"
$ws1.Range("F19").Value = 29
$ws1.Range("F20").Value = 560
$ws1.Range("F21").Value = 135
$ws1.Range("F24").Value = 2
$ws1.Range("F29").Value = 20

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 150
$ws2.Range("F12").Value = 60

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5710
$ws3.Range("F3").Value = 482
$ws3.Range("F4").Value = 470

# --- Sheet "全部类型" (aggregate of all the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 482
$ws4.Range("F4").Value = 470
$ws4.Range("F8").Value = 911
$ws4.Range("F10").Value = 464
$ws4.Range("F12").Value = 2191
$ws4.Range("F18").Value = 1088
$ws4.Range("F21").Value = 150
$ws4.Range("F24").Value = 668
$ws4.Range("F26").Value = 60
$ws4.Range("F27").Value = 1261
$ws4.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202409/UKoa3flf1726049294866.jpeg"
$ws4.Range("F28").Value = 29
$ws4.Range("F29").Value = 560
$ws4.Range("F30").Value = 135
$ws4.Range("F33").Value = 2
$ws4.Range("F49").Value = 20
